$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "20250523_084920"
$ws.Range("B4").Value = "2025-05-23 08:49:20"
$ws.Range("C4").Value = "Elie"
$ws.Range("D4").Value = "{'chantier': 'Atelier A', 'urgence': 'Normal', 'date_souhaitee': '2025-05-23', 'produits': {'102938475738883': {'produit': 'Parclose', 'quantite': 1, 'emplacement': 'Stockage'}}}"
$ws.Range("E4").Value = "f"
$ws.Range("F4").Value = "En attente"

# G4, H4, I4 stay blank, but still need a cell entry in the sheet XML.
# Touching a no-op style property materializes the cell without altering
# its formatting (keeps default style, avoids adding new style entries).
$ws.Range("G4").Font.Bold = $false
$ws.Range("H4").Font.Bold = $false
$ws.Range("I4").Font.Bold = $false
